# Update the "Predictor" column (C) values in the active worksheet so that
# several per-capita predictors are shown as natural-log transformed
# variables, matching the updated model specification.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value  = "ln(GDP [dollars per capita])"
$ws.Range("C3").Value  = "ln(ProMed Mentions [per capita])"
$ws.Range("C5").Value  = "ln(Tourism - Inbound [per capita])"
$ws.Range("C8").Value  = "ln(AB Exports [dollars per capita])"
$ws.Range("C9").Value  = "ln(Migrant Population [per capita])"
$ws.Range("C11").Value = "ln(Publication Bias Index [per capita])"
$ws.Range("C12").Value = "Livestock AB Consumption [kg per capita)"
$ws.Range("C14").Value = "ln(ProMed Mentions [per capita])"
$ws.Range("C15").Value = "ln(GDP [dollars per capita])"
$ws.Range("C16").Value = "ln(Publication Bias Index [per capita])"
$ws.Range("C17").Value = "ln(Population)"
